$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.409822305049147
$ws.Range("D2").Value = 0.02680039096887299
$ws.Range("E2").Value = 0.1732949958148993
$ws.Range("F2").Value = 0.7228684317892018
$ws.Range("G2").Value = 0.5640597308142645
$ws.Range("H2").Value = 0.685969125652548
$ws.Range("K2").Value = 1.155165454496455
$ws.Range("L2").Value = 0.1513846915122627
$ws.Range("N2").Value = 1.080579344301839
$ws.Range("O2").Value = 2.475141034919062
$ws.Range("C3").Value = 0.4005930589711113
$ws.Range("D3").Value = 0.02518600170974139
$ws.Range("E3").Value = 0.1689905923355326
$ws.Range("F3").Value = 0.7196873029621642
$ws.Range("G3").Value = 0.5623814605014701
$ws.Range("H3").Value = 0.6895562012097685
$ws.Range("K3").Value = 1.024494891603979
$ws.Range("L3").Value = 0.1472670776854841
$ws.Range("N3").Value = 1.078060015429742
$ws.Range("O3").Value = 2.478780715286888
$ws.Range("C4").Value = 0.3951462069826164
$ws.Range("D4").Value = 0.02418663106747942
$ws.Range("E4").Value = 0.1664433690756795
$ws.Range("F4").Value = 0.7182409409800457
$ws.Range("G4").Value = 0.5617985221254713
$ws.Range("H4").Value = 0.6921168147124206
$ws.Range("K4").Value = 0.94414387705379
$ws.Range("L4").Value = 0.1448246339277901
$ws.Range("N4").Value = 1.076913104609574
$ws.Range("O4").Value = 2.482702871458741
$ws.Range("C5").Value = 0.3929819528108283
$ws.Range("D5").Value = 0.02377736113879081
$ws.Range("E5").Value = 0.1654294373887062
$ws.Range("F5").Value = 0.7177788594551657
$ws.Range("G5").Value = 0.5616732716595578
$ws.Range("H5").Value = 0.6932503280300892
$ws.Range("K5").Value = 0.9113725382619009
$ws.Range("L5").Value = 0.1438508936184277
$ws.Range("N5").Value = 1.076546652857004
$ws.Range("O5").Value = 2.484724958417587
$ws.Range("C6").Value = 0.3926259272550681
$ws.Range("D6").Value = 0.02370928104290471
$ws.Range("E6").Value = 0.165262529878472
$ws.Range("F6").Value = 0.7177098172823051
$ws.Range("G6").Value = 0.5616592498511181
$ws.Range("H6").Value = 0.6934439850219007
$ws.Range("K6").Value = 0.9059292778800057
$ws.Range("L6").Value = 0.1436905079873725
$ws.Range("N6").Value = 1.07649190974881
$ws.Range("O6").Value = 2.485086304027817
$ws.Range("C7").Value = 0.3951167947521128
$ws.Range("D7").Value = 0.02418111964769309
$ws.Range("E7").Value = 0.1664295973043792
$ws.Range("F7").Value = 0.7182341938388745
$ws.Range("G7").Value = 0.5617963785582845
$ws.Range("H7").Value = 0.692131737109662
$ws.Range("K7").Value = 0.9437020203241389
$ws.Range("L7").Value = 0.1448114143626142
$ws.Range("N7").Value = 1.076907753410836
$ws.Range("O7").Value = 2.482728426881494
$ws.Range("C8").Value = 0.4065944535715005
$ws.Range("D8").Value = 0.02624545177954474
$ws.Range("E8").Value = 0.1717909846883288
$ws.Range("F8").Value = 0.7216662961505946
$ws.Range("G8").Value = 0.5633880369920092
$ws.Range("H8").Value = 0.6871316095711393
$ws.Range("K8").Value = 1.110136180542042
$ws.Range("L8").Value = 0.1499471333265276
$ws.Range("N8").Value = 1.079627870538175
$ws.Range("O8").Value = 2.476045456989937
$ws.Range("C9").Value = 0.4308455853532678
$ws.Range("D9").Value = 0.03022814260130247
$ws.Range("E9").Value = 0.1830639234750748
$ws.Range("F9").Value = 0.7324258551401854
$ws.Range("G9").Value = 0.570071775119473
$ws.Range("H9").Value = 0.6801690585705842
$ws.Range("K9").Value = 1.435489165050626
$ws.Range("L9").Value = 0.1606995269594762
$ws.Range("N9").Value = 1.088121828608777
$ws.Range("O9").Value = 2.476355842343594
$ws.Range("C10").Value = 0.4497259934568945
$ws.Range("D10").Value = 0.03311329647414851
$ws.Range("E10").Value = 0.1918100994285581
$ws.Range("F10").Value = 0.7427997487478706
$ws.Range("G10").Value = 0.5771717146191406
$ws.Range("H10").Value = 0.6767887420727874
$ws.Range("K10").Value = 1.673815348929224
$ws.Range("L10").Value = 0.169016474302623
$ws.Range("N10").Value = 1.096272485190113
$ws.Range("O10").Value = 2.48480461052759
$ws.Range("C11").Value = 0.4585462275133239
$ws.Range("D11").Value = 0.03441673972177739
$ws.Range("E11").Value = 0.1958900027451875
$ws.Range("F11").Value = 0.7480580097939793
$ws.Range("G11").Value = 0.5808809412073828
$ws.Range("H11").Value = 0.6756282085597576
$ws.Range("K11").Value = 1.782064933183221
$ws.Range("L11").Value = 0.1728911118761687
$ws.Range("N11").Value = 1.100392038037967
$ws.Range("O11").Value = 2.490442600034868
$ws.Range("C12").Value = 0.4619194639578552
$ws.Range("D12").Value = 0.03490899880754483
$ws.Range("E12").Value = 0.1974495114996344
$ws.Range("F12").Value = 0.7501268855562415
$ws.Range("G12").Value = 0.5823547508307456
$ws.Range("H12").Value = 0.6752430175499455
$ws.Range("K12").Value = 1.823030405006477
$ws.Range("L12").Value = 0.1743714674328629
$ws.Range("N12").Value = 1.102010930418658
$ws.Range("O12").Value = 2.492836294661799
$ws.Range("C13").Value = 0.4611915015330226
$ws.Range("D13").Value = 0.03480304145685409
$ws.Range("E13").Value = 0.1971129970964896
$ws.Range("F13").Value = 0.7496778584681749
$ws.Range("G13").Value = 0.5820342573970265
$ws.Range("H13").Value = 0.6753235606891366
$ws.Range("K13").Value = 1.814208965324781
$ws.Range("L13").Value = 0.1740520631507252
$ws.Range("N13").Value = 1.101659657866762
$ws.Range("O13").Value = 2.492309253170902
$ws.Range("C14").Value = 0.4588230806593856
$ws.Range("D14").Value = 0.03445726492270751
$ws.Range("E14").Value = 0.1960180132110452
$ws.Range("F14").Value = 0.7482266595786058
$ws.Range("G14").Value = 0.5810008039072869
$ws.Range("H14").Value = 0.6755954306332796
$ws.Range("K14").Value = 1.785435727933816
$ws.Range("L14").Value = 0.1730126388008131
$ws.Range("N14").Value = 1.100524046430152
$ws.Range("O14").Value = 2.490634342082416
$ws.Range("C15").Value = 0.4573766765340963
$ws.Range("D15").Value = 0.03424529320793823
$ws.Range("E15").Value = 0.195349196693634
$ws.Range("F15").Value = 0.7473478801791629
$ws.Range("G15").Value = 0.5803768043110153
$ws.Range("H15").Value = 0.6757690285473643
$ws.Range("K15").Value = 1.767807793609563
$ws.Range("L15").Value = 0.1723776691815431
$ws.Range("N15").Value = 1.099836114118048
$ws.Range("O15").Value = 2.489642121625792
$ws.Range("C16").Value = 0.4491542211434876
$ws.Range("D16").Value = 0.03302792930129783
$ws.Range("E16").Value = 0.191545503504571
$ws.Range("F16").Value = 0.7424669713736165
$ws.Range("G16").Value = 0.5769389779236462
$ws.Range("H16").Value = 0.6768721789651124
$ws.Range("K16").Value = 1.666737434625759
$ws.Range("L16").Value = 0.1687650929307694
$ws.Range("N16").Value = 1.096011522426281
$ws.Range("O16").Value = 2.484472322138942
$ws.Range("C17").Value = 0.4441692350106621
$ws.Range("D17").Value = 0.03227878351727043
$ws.Range("E17").Value = 0.1892379740895507
$ws.Range("F17").Value = 0.7396108895676718
$ws.Range("G17").Value = 0.5749529646595022
$ws.Range("H17").Value = 0.6776455605740068
$ws.Range("K17").Value = 1.604689768437481
$ws.Range("L17").Value = 0.1665722554726017
$ws.Range("N17").Value = 1.093770505543887
$ws.Range("O17").Value = 2.481760920507611
$ws.Range("C18").Value = 0.4413237917781032
$ws.Range("D18").Value = 0.0318470466361731
$ws.Range("E18").Value = 0.1879202729709135
$ws.Range("F18").Value = 0.7380188822768048
$ws.Range("G18").Value = 0.5738557726721325
$ws.Range("H18").Value = 0.6781258897539999
$ws.Range("K18").Value = 1.568986092007208
$ws.Range("L18").Value = 0.1653195792519995
$ws.Range("N18").Value = 1.092520308262976
$ws.Range("O18").Value = 2.480370271797796
$ws.Range("C19").Value = 0.4403641181670537
$ws.Range("D19").Value = 0.03170072316184047
$ws.Range("E19").Value = 0.1874757591332781
$ws.Range("F19").Value = 0.7374885649009428
$ws.Range("G19").Value = 0.5734920219736921
$ws.Range("H19").Value = 0.6782946170454238
$ws.Range("K19").Value = 1.556894851716322
$ws.Range("L19").Value = 0.1648969193126817
$ws.Range("N19").Value = 1.092103682268018
$ws.Range("O19").Value = 2.479928407176402
$ws.Range("C20").Value = 0.4446976407700731
$ws.Range("D20").Value = 0.03235861933987394
$ws.Range("E20").Value = 0.1894826284735416
$ws.Range("F20").Value = 0.7399096722554219
$ws.Range("G20").Value = 0.5751597084915971
$ws.Range("H20").Value = 0.67755955833006
$ws.Range("K20").Value = 1.611296470007801
$ws.Range("L20").Value = 0.1668047981738852
$ws.Range("N20").Value = 1.094005054647766
$ws.Range("O20").Value = 2.48203207091197
$ws.Range("C21").Value = 0.4595178427069868
$ws.Range("D21").Value = 0.0345588641010437
$ws.Range("E21").Value = 0.196339242110227
$ws.Range("F21").Value = 0.7486508022056597
$ws.Range("G21").Value = 0.581302473911947
$ws.Range("H21").Value = 0.6755141024473801
$ws.Range("K21").Value = 1.793887857979712
$ws.Range("L21").Value = 0.1733175870277819
$ws.Range("N21").Value = 1.100856006718175
$ws.Range("O21").Value = 2.491119277583721
$ws.Range("C22").Value = 0.4693972115210556
$ws.Range("D22").Value = 0.03598910578196524
$ws.Range("E22").Value = 0.2009051651290079
$ws.Range("F22").Value = 0.7548165257718153
$ws.Range("G22").Value = 0.5857206339665453
$ws.Range("H22").Value = 0.6744936696159556
$ws.Range("K22").Value = 1.913067502239244
$ws.Range("L22").Value = 0.1776505105521693
$ws.Range("N22").Value = 1.105676718694468
$ws.Range("O22").Value = 2.498566554189694
$ws.Range("C23").Value = 0.4641067253822087
$ws.Range("D23").Value = 0.0352264769525874
$ws.Range("E23").Value = 0.1984604992926364
$ws.Range("F23").Value = 0.7514842707270617
$ws.Range("G23").Value = 0.5833255705977081
$ws.Range("H23").Value = 0.6750093312748362
$ws.Range("K23").Value = 1.849473999570989
$ws.Range("L23").Value = 0.175330954340069
$ws.Range("N23").Value = 1.103072508750145
$ws.Range("O23").Value = 2.494453582596492
$ws.Range("C24").Value = 0.4444586846940695
$ws.Range("D24").Value = 0.03232252881073094
$ws.Range("E24").Value = 0.1893719924155874
$ws.Range("F24").Value = 0.73977443689121
$ws.Range("G24").Value = 0.5750661007186437
$ws.Range("H24").Value = 0.6775983287361811
$ws.Range("K24").Value = 1.608309678621254
$ws.Range("L24").Value = 0.1666996406529933
$ws.Range("N24").Value = 1.093898896003083
$ws.Range("O24").Value = 2.481908960011054
$ws.Range("C25").Value = 0.4240983358572805
$ws.Range("D25").Value = 0.02915783318123744
$ws.Range("E25").Value = 0.1799328978943464
$ws.Range("F25").Value = 0.729082455931227
$ws.Range("G25").Value = 0.5678804342313981
$ws.Range("H25").Value = 0.6817480691147466
$ws.Range("K25").Value = 1.347591145369449
$ws.Range("L25").Value = 0.1577175942206992
$ws.Range("N25").Value = 1.085487406993835
$ws.Range("O25").Value = 2.474831271567865
